$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value2 = 1088.75
$ws.Range("J41").Value2 = 1348.75
$ws.Range("L41").Value2 = 1348.75
$ws.Range("N41").Value2 = -2228.75

$ws.Range("H125").Value2 = 1138715.4
$ws.Range("I125").Value2 = 1138715.4
$ws.Range("K125").Value2 = 10248438.6
$ws.Range("M125").Value2 = -10245978.6

$ws.Range("H129").Value2 = 1927.3334
$ws.Range("I129").Value2 = 1498.75
$ws.Range("K129").Value2 = 4496.25
$ws.Range("M129").Value2 = 503.75

$ws.Range("H131").Value2 = 9746
$ws.Range("I131").Value2 = 9746
$ws.Range("K131").Value2 = 29238
$ws.Range("M131").Value2 = -24198

$ws.Range("H135").Value2 = 2264.3333
$ws.Range("J135").Value2 = 2733
$ws.Range("L135").Value2 = 24597
$ws.Range("N135").Value2 = -29667

$ws.Range("H137").Value2 = 4175509.8
$ws.Range("I137").Value2 = 7463.643
$ws.Range("J137").Value2 = 10010774
$ws.Range("K137").Value2 = 22390.929
$ws.Range("L137").Value2 = 30032322
$ws.Range("M137").Value2 = -19840.929
$ws.Range("N137").Value2 = -30037422

$ws.Range("H138").Value2 = 4548.0703
$ws.Range("I138").Value2 = 8352.643
$ws.Range("J138").Value2 = 3309.372
$ws.Range("K138").Value2 = 25057.929
$ws.Range("L138").Value2 = 9928.116
$ws.Range("M138").Value2 = -19917.929
$ws.Range("N138").Value2 = -20208.116

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 680.9
$ws.Range("I2").Value2 = 594.8
$ws.Range("J2").Value2 = 939.2
$ws.Range("K2").Value2 = 594.8
$ws.Range("L2").Value2 = 939.2
$ws.Range("M2").Value2 = -481.8
$ws.Range("N2").Value2 = -1165.2

$ws.Range("H31").Value2 = 4228.6665
$ws.Range("I31").Value2 = 4228.6665
$ws.Range("K31").Value2 = 4228.6665
$ws.Range("M31").Value2 = -3934.6665

$ws.Range("H46").Value2 = 7158.6
$ws.Range("I46").Value2 = 3989
$ws.Range("J46").Value2 = 7951
$ws.Range("K46").Value2 = 3989
$ws.Range("L46").Value2 = 7951
$ws.Range("N46").Value2 = -8589
$ws.Range("M46").Value2 = -3670

$ws.Range("H61").Value2 = 2440789.5
$ws.Range("I61").Value2 = 58400.45
$ws.Range("K61").Value2 = 58400.45
$ws.Range("M61").Value2 = -58188.45

$ws.Range("H88").Value2 = 1532.3334
$ws.Range("I88").Value2 = 1934.6666
$ws.Range("J88").Value2 = 1130
$ws.Range("K88").Value2 = 1934.6666
$ws.Range("L88").Value2 = 1130
$ws.Range("M88").Value2 = -1528.6666
$ws.Range("N88").Value2 = -1942

$ws.Range("H91").Value2 = 1532.3334
$ws.Range("I91").Value2 = 1934.6666
$ws.Range("J91").Value2 = 1130
$ws.Range("K91").Value2 = 1934.6666
$ws.Range("L91").Value2 = 1130
$ws.Range("M91").Value2 = -530.6666
$ws.Range("N91").Value2 = -3938

$ws.Range("H116").Value2 = 680.9
$ws.Range("I116").Value2 = 594.8
$ws.Range("J116").Value2 = 939.2
$ws.Range("K116").Value2 = 594.8
$ws.Range("L116").Value2 = 939.2
$ws.Range("M116").Value2 = 1699.2
$ws.Range("N116").Value2 = -5527.2

$ws.Range("H136").Value2 = 2440789.5
$ws.Range("I136").Value2 = 58400.45
$ws.Range("K136").Value2 = 175201.35
$ws.Range("M136").Value2 = -172651.35

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 680.9
$ws.Range("I3").Value2 = 594.8
$ws.Range("J3").Value2 = 939.2
$ws.Range("K3").Value2 = 594.8
$ws.Range("L3").Value2 = 939.2
$ws.Range("M3").Value2 = -480.8
$ws.Range("N3").Value2 = -1167.2

$ws.Range("H94").Value2 = 2955.6428
$ws.Range("I94").Value2 = 1614.0834
$ws.Range("K94").Value2 = 1614.0834
$ws.Range("M94").Value2 = -1163.0834

$ws.Range("H99").Value2 = 17302.125
$ws.Range("I99").Value2 = 31999.5
$ws.Range("K99").Value2 = 31999.5
$ws.Range("M99").Value2 = -30501.5

$ws.Range("H134").Value2 = 42859664
$ws.Range("I134").Value2 = 2186.8125
$ws.Range("J134").Value2 = 180003600
$ws.Range("K134").Value2 = 6560.4375
$ws.Range("L134").Value2 = 540010800
$ws.Range("M134").Value2 = -4025.4375
$ws.Range("N134").Value2 = -540015870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 2088.2424
$ws.Range("I58").Value2 = 2091.1538
$ws.Range("J58").Value2 = 2086.35
$ws.Range("K58").Value2 = 2091.1538
$ws.Range("L58").Value2 = 2086.35
$ws.Range("M58").Value2 = -1888.1538
$ws.Range("N58").Value2 = -2492.35

$ws.Range("H105").Value2 = 2341.25
$ws.Range("I105").Value2 = 1466.1111
$ws.Range("K105").Value2 = 1466.1111
$ws.Range("M105").Value2 = 280.8888999999999

$ws.Range("H107").Value2 = 986.8333
$ws.Range("I107").Value2 = 1132.3636
$ws.Range("J107").Value2 = 586.625
$ws.Range("K107").Value2 = 1132.3636
$ws.Range("L107").Value2 = 586.625
$ws.Range("M107").Value2 = 787.6364000000001
$ws.Range("N107").Value2 = -4426.625

$ws.Range("H132").Value2 = 16720392
$ws.Range("I132").Value2 = 75385.42999999999
$ws.Range("J132").Value2 = 55558740
$ws.Range("K132").Value2 = 226156.29
$ws.Range("L132").Value2 = 166676220
$ws.Range("M132").Value2 = -223626.29
$ws.Range("N132").Value2 = -166681280

$ws.Range("H134").Value2 = 673459.6
$ws.Range("J134").Value2 = 3124.25
$ws.Range("L134").Value2 = 9372.75
$ws.Range("N134").Value2 = -14442.75

$ws.Range("H136").Value2 = 2088.2424
$ws.Range("I136").Value2 = 2091.1538
$ws.Range("J136").Value2 = 2086.35
$ws.Range("K136").Value2 = 6273.4614
$ws.Range("L136").Value2 = 6259.049999999999
$ws.Range("M136").Value2 = -3723.4614
$ws.Range("N136").Value2 = -11359.05

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value2 = 180
$ws.Range("J21").Value2 = 0
$ws.Range("L21").Value2 = 0
$ws.Range("N21").ClearContents()

$ws.Range("H22").Value2 = 282838.22
$ws.Range("J22").Value2 = 500000
$ws.Range("L22").Value2 = 1500000
$ws.Range("N22").Value2 = -1500338

$ws.Range("H24").Value2 = 597.8570999999999
$ws.Range("J24").Value2 = 312.33334
$ws.Range("L24").Value2 = 937.0000200000001
$ws.Range("N24").Value2 = -1397.00002

$ws.Range("H27").Value2 = 282838.22
$ws.Range("J27").Value2 = 500000
$ws.Range("L27").Value2 = 1500000
$ws.Range("N27").Value2 = -1500204

$ws.Range("H52").Value2 = 1346.5
$ws.Range("J52").Value2 = 1346.5
$ws.Range("L52").Value2 = 4039.5
$ws.Range("N52").Value2 = -4571.5

$ws.Range("H132").Value2 = 2257.4546
$ws.Range("I132").Value2 = 2210.1428
$ws.Range("J132").Value2 = 2340.25
$ws.Range("K132").Value2 = 19891.2852
$ws.Range("L132").Value2 = 21062.25
$ws.Range("M132").Value2 = -17361.2852
$ws.Range("N132").Value2 = -26122.25

$ws.Range("H140").Value2 = 2713.5
$ws.Range("I140").Value2 = 1947
$ws.Range("K140").Value2 = 5841
$ws.Range("M140").Value2 = -661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value2 = 0
$ws.Range("J48").Value2 = 0
$ws.Range("L48").Value2 = 0
$ws.Range("N48").ClearContents()

$ws.Range("H95").Value2 = 59897
$ws.Range("J95").Value2 = 59897
$ws.Range("L95").Value2 = 59897
$ws.Range("N95").Value2 = -65389

$ws.Range("H113").Value2 = 1329.2
$ws.Range("I113").Value2 = 1329.2
$ws.Range("J113").Value2 = 0
$ws.Range("K113").Value2 = 1329.2
$ws.Range("L113").Value2 = 0
$ws.Range("M113").Value2 = 840.8
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 8599.950000000001
$ws.Range("I7").Value2 = 4331.3335
$ws.Range("J7").Value2 = 12092.454
$ws.Range("K7").Value2 = 4331.3335
$ws.Range("L7").Value2 = 12092.454
$ws.Range("M7").Value2 = -4219.3335
$ws.Range("N7").Value2 = -12316.454

$ws.Range("H93").Value2 = 1436.5
$ws.Range("J93").Value2 = 1883.3334
$ws.Range("L93").Value2 = 1883.3334
$ws.Range("N93").Value2 = -4379.3334

$ws.Range("H126").Value2 = 8599.950000000001
$ws.Range("I126").Value2 = 4331.3335
$ws.Range("J126").Value2 = 12092.454
$ws.Range("K126").Value2 = 12994.0005
$ws.Range("L126").Value2 = 36277.362
$ws.Range("M126").Value2 = -10524.0005
$ws.Range("N126").Value2 = -41217.362

$ws.Range("H136").Value2 = 64887.062
$ws.Range("I136").Value2 = 93549.91
$ws.Range("K136").Value2 = 280649.73
$ws.Range("M136").Value2 = -278099.73

$ws.Range("H140").Value2 = 71599.664
$ws.Range("J140").Value2 = 71599.664
$ws.Range("L140").Value2 = 71599.664
$ws.Range("N140").Value2 = -81959.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value2 = 1736.8889
$ws.Range("J96").Value2 = 869
$ws.Range("L96").Value2 = 869
$ws.Range("N96").Value2 = -3615
